$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data in column A for rows 2 and 3 (new shared strings "aaa", "bbb")
$ws.Range("A2").Value = "aaa"
$ws.Range("A3").Value = "bbb"

# Move the active selection to A4 (was B3)
[void]$ws.Range("A4").Select()
